# Actualización Automática de Datos (EA1, EA2 y EA3)
# Updates the "timestamp" column (H) for every data row (2-51) to the new
# refresh timestamp captured at ingestion time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-03-27 00:30:53"

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
